$wb = $excel.ActiveWorkbook

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 202.375
$ws.Range("I28").Value = 233.8
$ws.Range("J28").Value = 150
$ws.Range("K28").Value = 233.8
$ws.Range("L28").Value = 150
$ws.Range("M28").Value = 251.2
$ws.Range("N28").Value = -1120

# ALC row 57
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 21489.334
$ws.Range("J57").Value = 21489.334
$ws.Range("L57").Value = 64468.00199999999
$ws.Range("N57").Value = -65466.00199999999

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 31254616
$ws.Range("I135").Value = 1760.25
$ws.Range("J135").Value = 62507470
$ws.Range("K135").Value = 15842.25
$ws.Range("L135").Value = 562567230
$ws.Range("M135").Value = -13307.25
$ws.Range("N135").Value = -562572300

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1665.4736
$ws.Range("I137").Value = 1681.2273
$ws.Range("J137").Value = 1643.8125
$ws.Range("K137").Value = 5043.6819
$ws.Range("L137").Value = 4931.4375
$ws.Range("M137").Value = -2493.6819
$ws.Range("N137").Value = -10031.4375

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2087.5
$ws.Range("I141").Value = 1220.5883
$ws.Range("K141").Value = 3661.7649
$ws.Range("M141").Value = 1518.2351

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2111.84
$ws.Range("I45").Value = 1526.2222
$ws.Range("K45").Value = 1526.2222
$ws.Range("M45").Value = -1149.2222

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6511.875
$ws.Range("I61").Value = 8556
$ws.Range("J61").Value = 4922
$ws.Range("K61").Value = 8556
$ws.Range("L61").Value = 4922
$ws.Range("M61").Value = -8344
$ws.Range("N61").Value = -5346

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 65590.875
$ws.Range("I88").Value = 2003
$ws.Range("K88").Value = 2003
$ws.Range("M88").Value = -1597

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 65590.875
$ws.Range("I91").Value = 2003
$ws.Range("K91").Value = 2003
$ws.Range("M91").Value = -599

# ARM row 98
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H98").Value = 36333.168
$ws.Range("J98").Value = 36333.168
$ws.Range("L98").Value = 36333.168
$ws.Range("N98").Value = -42323.168

# ARM row 103
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H103").Value = 43120.668
$ws.Range("J103").Value = 43120.668
$ws.Range("L103").Value = 43120.668
$ws.Range("N103").Value = -45464.668

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6511.875
$ws.Range("I136").Value = 8556
$ws.Range("J136").Value = 4922
$ws.Range("K136").Value = 25668
$ws.Range("L136").Value = 14766
$ws.Range("M136").Value = -23118
$ws.Range("N136").Value = -19866

# BSM row 100
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 29845.428
$ws.Range("J100").Value = 29845.428
$ws.Range("L100").Value = 29845.428
$ws.Range("N100").Value = -32009.428

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2633266.5
$ws.Range("I105").Value = 1587.6364
$ws.Range("K105").Value = 1587.6364
$ws.Range("M105").Value = 159.3635999999999

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1081.2727
$ws.Range("I16").Value = 966
$ws.Range("K16").Value = 966
$ws.Range("M16").Value = -679

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2196.3076
$ws.Range("I31").Value = 975.03705
$ws.Range("J31").Value = 4944.1665
$ws.Range("K31").Value = 975.03705
$ws.Range("L31").Value = 4944.1665
$ws.Range("M31").Value = -680.03705
$ws.Range("N31").Value = -5534.1665

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2196.3076
$ws.Range("I34").Value = 975.03705
$ws.Range("J34").Value = 4944.1665
$ws.Range("K34").Value = 975.03705
$ws.Range("L34").Value = 4944.1665
$ws.Range("M34").Value = -773.03705
$ws.Range("N34").Value = -5348.1665

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 9616539
$ws.Range("I105").Value = 12500740
$ws.Range("J105").Value = 2533
$ws.Range("K105").Value = 12500740
$ws.Range("L105").Value = 2533
$ws.Range("M105").Value = -12498993
$ws.Range("N105").Value = -6027

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1081.2727
$ws.Range("I113").Value = 966
$ws.Range("K113").Value = 966
$ws.Range("M113").Value = 1204

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2986.6365
$ws.Range("I132").Value = 868.46155
$ws.Range("J132").Value = 6046.222
$ws.Range("K132").Value = 2605.38465
$ws.Range("L132").Value = 18138.666
$ws.Range("M132").Value = -75.38464999999997
$ws.Range("N132").Value = -23198.666

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 805.3333
$ws.Range("I134").Value = 805.3333
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 2415.9999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 119.0001000000002
$ws.Range("N134").ClearContents()

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1128.6364
$ws.Range("I68").Value = 733.3333
$ws.Range("K68").Value = 2199.9999
$ws.Range("M68").Value = -1388.9999

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1128.6364
$ws.Range("I71").Value = 733.3333
$ws.Range("K71").Value = 6599.9997
$ws.Range("M71").Value = -2543.9997

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 3022.85
$ws.Range("I107").Value = 4586.75
$ws.Range("J107").Value = 677
$ws.Range("K107").Value = 13760.25
$ws.Range("L107").Value = 2031
$ws.Range("M107").Value = -11840.25
$ws.Range("N107").Value = -5871

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 416.66666
$ws.Range("I113").Value = 431.66666
$ws.Range("J113").Value = 386.66666
$ws.Range("K113").Value = 1294.99998
$ws.Range("L113").Value = 1159.99998
$ws.Range("M113").Value = 875.0000199999999
$ws.Range("N113").Value = -5499.999980000001

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1101.6
$ws.Range("J122").Value = 1251
$ws.Range("L122").Value = 11259
$ws.Range("N122").Value = -16159

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 806.52
$ws.Range("J131").Value = 828.54736
$ws.Range("L131").Value = 2485.64208
$ws.Range("N131").Value = -12565.64208

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3220.88
$ws.Range("I80").Value = 2685.3845
$ws.Range("J80").Value = 3801
$ws.Range("K80").Value = 2685.3845
$ws.Range("L80").Value = 3801
$ws.Range("M80").Value = -1687.3845
$ws.Range("N80").Value = -5797

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3220.88
$ws.Range("I83").Value = 2685.3845
$ws.Range("J83").Value = 3801
$ws.Range("K83").Value = 13426.9225
$ws.Range("L83").Value = 19005
$ws.Range("M83").Value = -8434.922500000001
$ws.Range("N83").Value = -28989

# GSM row 98
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 21597
$ws.Range("J98").Value = 21597
$ws.Range("L98").Value = 21597
$ws.Range("N98").Value = -27587

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2135.8
$ws.Range("I122").Value = 2154.75
$ws.Range("J122").Value = 2060
$ws.Range("K122").Value = 6464.25
$ws.Range("L122").Value = 6180
$ws.Range("M122").Value = -4014.25
$ws.Range("N122").Value = -11080

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3073
$ws.Range("I122").Value = 2649
$ws.Range("J122").Value = 3436.4285
$ws.Range("K122").Value = 7947
$ws.Range("L122").Value = 10309.2855
$ws.Range("M122").Value = -5497
$ws.Range("N122").Value = -15209.2855

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1569.1111
$ws.Range("I132").Value = 1017.9
$ws.Range("J132").Value = 3144
$ws.Range("K132").Value = 3053.7
$ws.Range("L132").Value = 9432
$ws.Range("M132").Value = -523.6999999999998
$ws.Range("N132").Value = -14492

# LTW row 140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 46452.668
$ws.Range("J140").Value = 46452.668
$ws.Range("L140").Value = 46452.668
$ws.Range("N140").Value = -56812.668

# WVR row 98
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 49333.332
$ws.Range("J98").Value = 49333.332
$ws.Range("L98").Value = 49333.332
$ws.Range("N98").Value = -55323.332

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1423.6
$ws.Range("I132").Value = 1259.8462
$ws.Range("J132").Value = 1727.7142
$ws.Range("K132").Value = 3779.5386
$ws.Range("L132").Value = 5183.142599999999
$ws.Range("M132").Value = -1249.5386
$ws.Range("N132").Value = -10243.1426

# WVR row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 44286.25
$ws.Range("J135").Value = 44286.25
$ws.Range("L135").Value = 44286.25
$ws.Range("N135").Value = -54426.25
